# Update the "取得日時" (acquired timestamp) column on the ランサーズ sheet
# from 2026-01-17 18:26:13 to 2026-01-17 18:33:52 for rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = "2026-01-17 18:33:52"
}
